$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.921.20"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "2.582.74"
$ws.Range("E3").Value = "  -4.99%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.104"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "3.040.54"
$ws.Range("E13").Value = "  -4.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("D15").Value = "61.807.62"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("D17").Value = "2.589.89"
$ws.Range("E17").Value = "  -4.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.496"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D28").Value = "0.0₃0836"
$ws.Range("E28").Value = "  -6.37%  "
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.15%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("E36").Value = "  -4.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "338.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.888"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.56%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.134.61"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.607"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0547"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.86%  "
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("E51").Value = "  -2.61%  "
